$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-24 22:48:55"
$wsZhCn.Range("H4").Value = "2016-03-24 22:49:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-24 22:49:00"
$wsDeDe.Range("H4").Value = "2016-03-24 22:49:33"
